{"js": "// The template's second paragraph holds a Word FIELD whose instrText\n// spells out \" m:'contents.txt'.fromConfluenceURI() \". The edit turns\n// that field into plain literal text \"{m:'contents.txt'.fromConfluenceURI()}\"\n// (the M2Doc template syntax), keeping the _GoBack bookmark in place,\n// and removing the field begin/end characters entirely.\n\nconst body = context.document.body;\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the paragraph that currently only contains the field (it is\n// the second paragraph in this template, right after the \"query :\"\n// intro paragraph).\nconst fieldParagraph = paragraphs.items[1];\n\n// Find and delete the field itself (removes the begin fldChar, all the\n// instrText runs, and the end fldChar).\nconst fields = body.fields;\nfields.load(\"items\");\nawait context.sync();\n\nif (fields.items.length > 0) {\n    fields.items[0].delete();\n    await context.sync();\n}\n\n// Re-insert the same characters as plain text runs (one run per former\n// instrText run, to mirror how the template is written), wrapping the\n// whole thing in curly braces instead of field delimiters, and keeping\n// the bookmark exactly where it was (between \"Confluence\" and \"URI\").\nconst runsXml =\n    \"<w:r><w:t>{</w:t></w:r>\" +\n    \"<w:r><w:t>m</w:t></w:r>\" +\n    \"<w:r><w:t>:</w:t></w:r>\" +\n    \"<w:r><w:t>'</w:t></w:r>\" +\n    \"<w:r><w:t>contents.txt</w:t></w:r>\" +\n    \"<w:r><w:t>'</w:t></w:r>\" +\n    \"<w:r><w:t>.from</w:t></w:r>\" +\n    \"<w:r><w:t>Confluence</w:t></w:r>\" +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    \"<w:r><w:t>URI</w:t></w:r>\" +\n    \"<w:r><w:t>()</w:t></w:r>\" +\n    '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>';\n\nconst ooxml =\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p>\" + runsXml + \"</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\nfieldParagraph.insertOoxml(ooxml, \"Start\");\nawait context.sync();\n", "ps1": "# The template's second paragraph holds a Word FIELD whose instrText\n# spells out \" m:'contents.txt'.fromConfluenceURI() \". This edit turns\n# that field into plain literal text \"{m:'contents.txt'.fromConfluenceURI()}\"\n# (the M2Doc template syntax), keeping the _GoBack bookmark in place,\n# and removing the field begin/end characters entirely.\n\n$d = $word.ActiveDocument\n\n# The field-bearing paragraph is the 2nd paragraph of the document (right\n# after the \"A simple demonstration of a query :\" intro paragraph).\n$fieldParagraph = $d.Paragraphs.Item(2)\n\n# Delete the field itself -- this removes the begin fldChar, all of the\n# instrText runs, and the end fldChar from the paragraph.\n$f = $d.Fields.Item(1)\n$f.Delete()\n\n# Collapsed insertion point at the (now empty) start of that paragraph,\n# so we insert content without touching the paragraph mark / <w:p> attrs.\n$insertPoint = $d.Range($fieldParagraph.Range.Start, $fieldParagraph.Range.Start)\n\n# Re-insert the same characters as plain text runs (one run per former\n# instrText run, mirroring how the template was written), wrapping the\n# whole thing in curly braces instead of field delimiters, and keeping\n# the bookmark exactly where it was (between \"Confluence\" and \"URI\").\n$runsXml = \"<w:r><w:t>{</w:t></w:r>\" + `\n    \"<w:r><w:t>m</w:t></w:r>\" + `\n    \"<w:r><w:t>:</w:t></w:r>\" + `\n    \"<w:r><w:t>'</w:t></w:r>\" + `\n    \"<w:r><w:t>contents.txt</w:t></w:r>\" + `\n    \"<w:r><w:t>'</w:t></w:r>\" + `\n    \"<w:r><w:t>.from</w:t></w:r>\" + `\n    \"<w:r><w:t>Confluence</w:t></w:r>\" + `\n    \"<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>\" + `\n    \"<w:r><w:t>URI</w:t></w:r>\" + `\n    \"<w:r><w:t>()</w:t></w:r>\" + `\n    \"<w:r><w:t xml:space='preserve'>}</w:t></w:r>\"\n\n$ooxml = \"<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'>\" + `\n    \"<pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'>\" + `\n    \"<pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>\" + `\n    \"<w:body><w:p>$runsXml</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\"\n\n$insertPoint.InsertXML($ooxml)\n"}
